# Atualização de bases das ligas, do dia: 21-04-2024 às 14:32
#
# The sheet currently ends at row 94 (id=92, match 7802875, York United FC vs
# Forge FC, played 2024-04-23). Two newer (earlier-dated) match results are
# inserted ahead of it, and the old row 94 is pushed down to row 96 (its id
# is renumbered from 92 to 94 to keep the running sequence).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: preserve the current last row (id=92) by copying it down to
# row 96, cell by cell (only the columns that actually hold data), so
# styles (bold/border id style, date number format) come along without
# minting new style entries and without materialising stray blank cells
# in columns that were never populated for this row (H, I, J, AB, AC).
# ---------------------------------------------------------------------
$oldRowCols = @("A","B","C","D","E","F","G","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")
foreach ($col in $oldRowCols) {
    $srcAddr = $col + "94"
    $dstAddr = $col + "96"
    $ws.Range($srcAddr).Copy($ws.Range($dstAddr))
}
# that match's running id moves from 92 -> 94
$ws.Range("A96").Value2 = 94

# ---------------------------------------------------------------------
# Step 2: overwrite row 94 with the first newly-added match
# (id=92, match 7802935, Pacific FC CA vs Valour FC, 2024-04-21, 2-0 H)
# ---------------------------------------------------------------------
$ws.Range("A93").Copy($ws.Range("A94"))
$ws.Range("A94").Value2 = 92

$ws.Range("B94").Value2 = 7802935
$ws.Range("C94").Value2 = "Canada Premier League"
$ws.Range("D94").Value2 = "Canada Premier League"

$ws.Range("E93").Copy($ws.Range("E94"))
$ws.Range("E94").Value2 = 45401.95833333334

$ws.Range("F94").Value2 = "Pacific FC CA"
$ws.Range("G94").Value2 = "Valour FC"
$ws.Range("H94").Value2 = 2
$ws.Range("I94").Value2 = 0
$ws.Range("J94").Value2 = "H"
$ws.Range("K94").Value2 = 1.727
$ws.Range("L94").Value2 = 3.5
$ws.Range("M94").Value2 = 4
$ws.Range("N94").Value2 = 1.615
$ws.Range("O94").Value2 = 4
$ws.Range("P94").Value2 = 4.2
$ws.Range("Q94").Value2 = -0.75
$ws.Range("R94").Value2 = 1.9
$ws.Range("S94").Value2 = 1.9
$ws.Range("T94").Value2 = 2.5
$ws.Range("U94").Value2 = 1.95
$ws.Range("V94").Value2 = 1.75
$ws.Range("W94").Value2 = 0.615
$ws.Range("X94").Value2 = -1
$ws.Range("Y94").Value2 = -1
$ws.Range("Z94").Value2 = 0.8999999999999999
$ws.Range("AA94").Value2 = -1
$ws.Range("AB94").Value2 = -1
$ws.Range("AC94").Value2 = 0.75

# ---------------------------------------------------------------------
# Step 3: write row 95 with the second newly-added match
# (id=93, match 7802936, Atletico Ottawa vs Cavalry FC, 2024-04-22, 1-1 D)
# ---------------------------------------------------------------------
$ws.Range("A93").Copy($ws.Range("A95"))
$ws.Range("A95").Value2 = 93

$ws.Range("B95").Value2 = 7802936
$ws.Range("C95").Value2 = "Canada Premier League"
$ws.Range("D95").Value2 = "Canada Premier League"

$ws.Range("E93").Copy($ws.Range("E95"))
$ws.Range("E95").Value2 = 45402.625

$ws.Range("F95").Value2 = "Atletico Ottawa"
$ws.Range("G95").Value2 = "Cavalry FC"
$ws.Range("H95").Value2 = 1
$ws.Range("I95").Value2 = 1
$ws.Range("J95").Value2 = "D"
$ws.Range("K95").Value2 = 3.1
$ws.Range("L95").Value2 = 3.4
$ws.Range("M95").Value2 = 2
$ws.Range("N95").Value2 = 2.875
$ws.Range("O95").Value2 = 3.1
$ws.Range("P95").Value2 = 2.25
$ws.Range("Q95").Value2 = 0.25
$ws.Range("R95").Value2 = 1.775
$ws.Range("S95").Value2 = 2.025
$ws.Range("T95").Value2 = 2.25
$ws.Range("U95").Value2 = 1.95
$ws.Range("V95").Value2 = 1.85
$ws.Range("W95").Value2 = -1
$ws.Range("X95").Value2 = 2.1
$ws.Range("Y95").Value2 = -1
$ws.Range("Z95").Value2 = 0.3875
$ws.Range("AA95").Value2 = -0.5
$ws.Range("AB95").Value2 = -0.5
$ws.Range("AC95").Value2 = 0.425

Write-Output "done"
